{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst texts = paragraphs.items.map((p) => p.text);\n\n// Locate the \"Ver no Jupiter Salvar em pdf Salvar em docx\" paragraph and the\n// \"\u00a9 2020 ... Contact: luizeleno@usp.br ...\" paragraph right after it.\nconst jupiterIndex = texts.findIndex((t) => t.indexOf(\"Ver no Jupiter\") === 0);\nconst copyrightIndex = texts.findIndex(\n  (t) => t.indexOf(\"Contact: luizeleno@usp.br\") !== -1\n);\n\nif (jupiterIndex !== -1 && copyrightIndex !== -1) {\n  // Also remove the blank separator paragraph immediately preceding\n  // \"Ver no Jupiter...\".\n  const blankIndex = jupiterIndex - 1;\n\n  // Delete from the bottom up so earlier indices stay valid.\n  paragraphs.items[copyrightIndex].delete();\n  paragraphs.items[jupiterIndex].delete();\n  if (blankIndex >= 0) {\n    paragraphs.items[blankIndex].delete();\n  }\n\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"Ver no Jupiter Salvar em pdf Salvar em docx\" paragraph and the\n# \"(c) 2020 ... Contact: luizeleno@usp.br ...\" paragraph that immediately\n# follows it. These two paragraphs (plus the blank separator paragraph right\n# before the first one) are being removed from the bottom of the page.\n$jupiterIndex = $null\n$copyrightIndex = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $paraText = $d.Paragraphs.Item($i).Range.Text\n    if ($jupiterIndex -eq $null -and $paraText -like \"Ver no Jupiter*\") {\n        $jupiterIndex = $i\n    }\n    if ($copyrightIndex -eq $null -and $paraText -like \"*Contact: luizeleno@usp.br*\") {\n        $copyrightIndex = $i\n    }\n}\n\nif ($jupiterIndex -ne $null -and $copyrightIndex -ne $null) {\n    # Also remove the blank paragraph immediately preceding \"Ver no Jupiter...\".\n    $blankIndex = $jupiterIndex - 1\n\n    $startRange = $d.Paragraphs.Item($blankIndex).Range\n    $endRange = $d.Paragraphs.Item($copyrightIndex).Range\n\n    $deleteRange = $d.Range($startRange.Start, $endRange.End)\n    $deleteRange.Delete()\n}\n"}
